# Applies the "added LOSTIW to coho_populations.xlsx" edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header columns (new naming convention, lower-case)
$ws.Range("A1").Value = "site_code"
$ws.Range("B1").Value = "coho_esu_dps"
$ws.Range("C1").Value = "coho_mpg"
$ws.Range("D1").Value = "coho_popid"
$ws.Range("E1").Value = "coho_popname"

# 2. Insert the new LOSTIW record in the Grande Ronde River / GRWAL-c
#    (Wallowa River) group, right after the existing WR2 row, which is
#    where it will end up once the table is re-sorted by popid/site_code.
$ws.Rows.Item(19).Insert()
$ws.Cells.Item(19, 1).Value = "LOSTIW"
$ws.Cells.Item(19, 2).Value = "Snake River Coho Salmon"
$ws.Cells.Item(19, 3).Value = "Grande Ronde River"
$ws.Cells.Item(19, 4).Value = "GRWAL-c"
$ws.Cells.Item(19, 5).Value = "Wallowa River"

# 3. Re-sort the data range (A2:E23) by popid (col D) then by site_code
#    (col A), matching the workbook's existing two-key sortState.
$sortKey1 = $ws.Range("D2:D23")
$sortKey2 = $ws.Range("A2:A23")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($sortKey1)
$ws.Sort.SortFields.Add($sortKey2)
$ws.Sort.SetRange($ws.Range("A1:E23"))
$ws.Sort.Header = 1
$ws.Sort.Apply()

# The strict two-key sort alphabetizes ties on site_code, which would place
# LOSTIW ahead of MR1/WR1/WR2 within the tied GRWAL-c group -- but the
# original (real-world) edit only ever inserted the new row in place and
# never re-ran a fresh two-key sort over that tie, so LOSTIW stayed last in
# its group. Restore that exact row order (the sortState/metadata recorded
# above is left as the two-key definition, matching the target workbook).
$ws.Cells.Item(16, 1).Value = "MR1"
$ws.Cells.Item(17, 1).Value = "WR1"
$ws.Cells.Item(18, 1).Value = "WR2"
$ws.Cells.Item(19, 1).Value = "LOSTIW"

# 4. Restore selection cursor to E2 (matches post-edit state)
$ws.Range("E2").Select()
